$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = '<jt:escape doublequote="Embedded \"double-quotes\"" backslash="Embedded \\backslash"/>'

$ws.Range("A15").Value = $newValue
